$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 61.5
$ws.Range("I8").Value = 61.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 184.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -45.5
$ws.Range("N8").ClearContents()
$ws.Range("H17").Value = 5711100
$ws.Range("J17").Value = 5711100
$ws.Range("L17").Value = 17133300
$ws.Range("N17").Value = -17133636
$ws.Range("H138").Value = 2331.75
$ws.Range("I138").Value = 1544.174
$ws.Range("J138").Value = 2610.4307
$ws.Range("K138").Value = 4632.522
$ws.Range("L138").Value = 7831.2921
$ws.Range("M138").Value = 507.4780000000001
$ws.Range("N138").Value = -18111.2921

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 28864.75
$ws.Range("J37").Value = 28864.75
$ws.Range("L37").Value = 28864.75
$ws.Range("N37").Value = -29410.75
$ws.Range("H61").Value = 308682.66
$ws.Range("I61").Value = 6268.6
$ws.Range("J61").Value = 1253726.6
$ws.Range("K61").Value = 6268.6
$ws.Range("L61").Value = 1253726.6
$ws.Range("M61").Value = -6056.6
$ws.Range("N61").Value = -1254150.6
$ws.Range("H75").Value = 42000
$ws.Range("J75").Value = 42000
$ws.Range("L75").Value = 42000
$ws.Range("N75").Value = -43748
$ws.Range("H78").Value = 42000
$ws.Range("J78").Value = 42000
$ws.Range("L78").Value = 126000
$ws.Range("N78").Value = -134736
$ws.Range("H136").Value = 308682.66
$ws.Range("I136").Value = 6268.6
$ws.Range("J136").Value = 1253726.6
$ws.Range("K136").Value = 18805.8
$ws.Range("L136").Value = 3761179.8
$ws.Range("M136").Value = -16255.8
$ws.Range("N136").Value = -3766279.8

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1804.6
$ws.Range("I8").Value = 1804.6
$ws.Range("K8").Value = 1804.6
$ws.Range("M8").Value = -1664.6
$ws.Range("H11").Value = 766.25
$ws.Range("J11").Value = 1897.5
$ws.Range("L11").Value = 1897.5
$ws.Range("N11").Value = -2177.5
$ws.Range("H94").Value = 2052.7693
$ws.Range("I94").Value = 1415.1428
$ws.Range("J94").Value = 2796.6667
$ws.Range("K94").Value = 1415.1428
$ws.Range("L94").Value = 2796.6667
$ws.Range("M94").Value = -964.1428000000001
$ws.Range("N94").Value = -3698.6667

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2664.125
$ws.Range("I16").Value = 2566.6667
$ws.Range("J16").Value = 2722.6
$ws.Range("K16").Value = 2566.6667
$ws.Range("L16").Value = 2722.6
$ws.Range("M16").Value = -2279.6667
$ws.Range("N16").Value = -3296.6
$ws.Range("H113").Value = 2664.125
$ws.Range("I113").Value = 2566.6667
$ws.Range("J113").Value = 2722.6
$ws.Range("K113").Value = 2566.6667
$ws.Range("L113").Value = 2722.6
$ws.Range("M113").Value = -396.6667000000002
$ws.Range("N113").Value = -7062.6

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 342.33334
$ws.Range("J6").Value = 483.33334
$ws.Range("L6").Value = 1450.00002
$ws.Range("N6").Value = -1676.00002
$ws.Range("H23").Value = 5882469
$ws.Range("I23").Value = 16666745
$ws.Range("J23").Value = 136.27272
$ws.Range("K23").Value = 50000235
$ws.Range("L23").Value = 408.81816
$ws.Range("M23").Value = -50000000
$ws.Range("N23").Value = -878.81816
$ws.Range("H122").Value = 4110.8066
$ws.Range("I122").Value = 580.94446
$ws.Range("J122").Value = 8998.308000000001
$ws.Range("K122").Value = 5228.50014
$ws.Range("L122").Value = 80984.77200000001
$ws.Range("M122").Value = -2778.50014
$ws.Range("N122").Value = -85884.77200000001
$ws.Range("H139").Value = 3940.8125
$ws.Range("I139").Value = 4790.36
$ws.Range("J139").Value = 3017.3914
$ws.Range("K139").Value = 14371.08
$ws.Range("L139").Value = 9052.174199999999
$ws.Range("M139").Value = -9231.079999999998
$ws.Range("N139").Value = -19332.1742

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 590060
$ws.Range("I7").Value = 666766.7
$ws.Range("J7").Value = 475000
$ws.Range("K7").Value = 666766.7
$ws.Range("L7").Value = 475000
$ws.Range("M7").Value = -666654.7
$ws.Range("N7").Value = -475224
$ws.Range("H8").Value = 590060
$ws.Range("I8").Value = 666766.7
$ws.Range("J8").Value = 475000
$ws.Range("K8").Value = 666766.7
$ws.Range("L8").Value = 475000
$ws.Range("M8").Value = -666627.7
$ws.Range("N8").Value = -475278
$ws.Range("H34").Value = 39000
$ws.Range("I34").Value = 39000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 39000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -38732
$ws.Range("N34").ClearContents()
$ws.Range("H76").Value = 39000
$ws.Range("I76").Value = 39000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 39000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -38685
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 39000
$ws.Range("I79").Value = 39000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 39000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -37908
$ws.Range("N79").ClearContents()
$ws.Range("H113").Value = 58825430
$ws.Range("I113").Value = 100001270
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 100001270
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = -99999100
$ws.Range("N113").Value = -7140

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 337999.66
$ws.Range("I18").Value = 503499.5
$ws.Range("J18").Value = 7000
$ws.Range("K18").Value = 503499.5
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = -503327.5
$ws.Range("N18").Value = -7344
$ws.Range("H61").Value = 1724.5454
$ws.Range("I61").Value = 1617
$ws.Range("K61").Value = 1617
$ws.Range("M61").Value = -1415
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H75").Value = 42200
$ws.Range("J75").Value = 42200
$ws.Range("L75").Value = 42200
$ws.Range("N75").Value = -44072
$ws.Range("H76").Value = 15798
$ws.Range("J76").Value = 15798
$ws.Range("L76").Value = 15798
$ws.Range("N76").Value = -16474
$ws.Range("H78").Value = 42200
$ws.Range("J78").Value = 42200
$ws.Range("L78").Value = 126600
$ws.Range("N78").Value = -135960
$ws.Range("H79").Value = 15798
$ws.Range("J79").Value = 15798
$ws.Range("L79").Value = 15798
$ws.Range("N79").Value = -18138
$ws.Range("H113").Value = 1724.5454
$ws.Range("I113").Value = 1617
$ws.Range("K113").Value = 1617
$ws.Range("M113").Value = 553

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 450
$ws.Range("I7").Value = 450
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 450
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -337
$ws.Range("N7").ClearContents()
$ws.Range("H100").Value = 451.25
$ws.Range("I100").Value = 301
$ws.Range("J100").Value = 601.5
$ws.Range("K100").Value = 602
$ws.Range("L100").Value = 1203
$ws.Range("M100").Value = -61
$ws.Range("N100").Value = -2285
$ws.Range("H113").Value = 2251.6875
$ws.Range("I113").Value = 2086.0454
$ws.Range("J113").Value = 2616.1
$ws.Range("K113").Value = 6258.1362
$ws.Range("L113").Value = 7848.299999999999
$ws.Range("M113").Value = -4088.1362
$ws.Range("N113").Value = -12188.3
$ws.Range("H132").Value = 1823.7567
$ws.Range("I132").Value = 1370.5
$ws.Range("J132").Value = 2488.5334
$ws.Range("K132").Value = 4111.5
$ws.Range("L132").Value = 7465.600199999999
$ws.Range("M132").Value = -1581.5
$ws.Range("N132").Value = -12525.6002
